$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[Cláudio-Tecnologia da Soldagem, Cleidson-Automação Industrial, Leandro-Sistemas de Refrigeração, Paulo Rob.-CAM]"
$ws.Range("C18").Value = "[Paulo Rob.-CAM, Guilherme-Eletropneumática, Cláudio-Tecnologia da Soldagem, Leandro-Sistemas de Refrigeração]"
$ws.Range("D18").Value = "[Allan Cupertino-Instalções Elétricas, Emerson-Eletrônica Básica]"
$ws.Range("E18").Value = "Allan Cupertino-Máquinas Elétri"
$ws.Range("F18").Value = "[-, Allan Cupertino-Lab. De Máquinas elétricas]"

$ws.Range("B19").Value = "[Guilherme-Eletro, Cleidson-Automação Industrial, Guilherme-Eletropneumática, Paulo Rob.-CAM]"
$ws.Range("C19").Value = "Andre B.-Circuitos Elétrico"
$ws.Range("D19").Value = "[Emerson-Eletrônica Básica, Emerson-Eletrônica Básica]"
$ws.Range("E19").Value = "[Weslei-CAD, Weslei-CAD]"
$ws.Range("F19").Value = "[-, João Paulo-Lab. Circuitos Elétricos]"

$ws.Range("B20").Value = "[Guilherme-Eletro, Cleidson-Automação Industrial, Leandro-Sistemas de Refrigeração, Paulo Rob.-CAM]"
$ws.Range("C20").Value = "[Guilherme-Eletro, Guilherme-Eletropneumática, Cláudio-Tecnologia da Soldagem, Leandro-Sistemas de Refrigeração]"
$ws.Range("D20").Value = "[Emerson-Eletrônica Básica, Allan Cupertino-Instalções Elétricas]"
$ws.Range("E20").Value = "[-, Weslei-CAD]"
$ws.Range("F20").Value = "[Weslei-CAD, Allan Cupertino-Instalções Elétricas]"

$ws.Range("B21").Value = "[Guilherme-Eletro, Cleidson-Automação Industrial, Cláudio-Tecnologia da Soldagem, Guilherme-Eletropneumática]"
$ws.Range("C21").Value = "Andre B.-Circuitos Elétrico"
$ws.Range("D21").Value = "[Allan Cupertino-Lab. De Máquinas elétricas, João Paulo-Lab. Circuitos Elétricos]"
$ws.Range("E21").Value = "Allan Cupertino-Máquinas Elétri"
$ws.Range("F21").Value = "[-, Allan Cupertino-Instalções Elétricas]"
